$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scraper added two new columns ("height" and "weight") between the
# existing "fumbles" (D) and "fantasy points" (E) columns, pushing the old
# "fantasy points" column from E to G. Replicate that: give F and G the
# same header formatting as the existing headers (bold, bordered, centered)
# by copying E1's format, then fill in the new values.
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("E1").Copy($ws.Range("G1"))

# Shift each row's existing "fantasy points" value (currently in column E)
# over to the new column G, then populate the new height/weight columns.
for ($r = 2; $r -le 17; $r++) {
    $fantasyPoints = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 7).Value = $fantasyPoints
    $ws.Cells.Item($r, 5).Value = 6.416666666666667
    $ws.Cells.Item($r, 6).Value = 237
}

# Now that the old E-column values have been copied over to G, relabel the
# header row: E=height, F=weight, G=fantasy points.
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"
$ws.Range("G1").Value = "fantasy points"
